$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.167.13'
$ws.Range('E2').Value = '  -2.20%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.576.39'
$ws.Range('E3').Value = '  -1.74%  '

$ws.Range('E4').Value = '  -0.54%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '208.74'
$ws.Range('E5').Value = '  -1.55%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  -3.00%  '

$ws.Range('E7').Value = '  -0.44%  '

$ws.Range('E8').Value = '  -1.66%  '

$ws.Range('E9').Value = '  -0.93%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.29%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  -0.57%  '

$ws.Range('E12').Value = '  -1.67%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.06'
$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.569.04'
$ws.Range('E14').Value = '  -1.45%  '

$ws.Range('E15').Value = '  -2.08%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.40'
$ws.Range('E16').Value = '  -1.07%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.160.16'
$ws.Range('E17').Value = '  -2.15%  '

$ws.Range('E18').Value = '  -1.95%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.25'
$ws.Range('E19').Value = '  +1.17%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '208.89'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('E21').Value = '  -0.53%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -1.23%  '

$ws.Range('E23').Value = '  -2.56%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.83'
$ws.Range('E24').Value = '  -2.28%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '143.87'
$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('E26').Value = '  -0.55%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.98'
$ws.Range('E27').Value = '  -1.57%  '

$ws.Range('E28').Value = '  -1.80%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.21'
$ws.Range('E29').Value = '  -1.06%  '

$ws.Range('E30').Value = '  -0.33%  '

$ws.Range('E31').Value = '  -1.58%  '

$ws.Range('E32').Value = '  -1.92%  '

$ws.Range('E33').Value = '  +1.14%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.279.49'
$ws.Range('E34').Value = '  -0.79%  '

$ws.Range('E35').Value = '  -1.66%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.608'
$ws.Range('E36').Value = '  +3.52%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.47'
$ws.Range('E37').Value = '  -1.43%  '

$ws.Range('E38').Value = '  -2.13%  '

$ws.Range('E39').Value = '  -9.60%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  -1.92%  '

$ws.Range('E41').Value = '  -0.45%  '

$ws.Range('E42').Value = '  +2.71%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.765'
$ws.Range('E43').Value = '  -1.83%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.13'
$ws.Range('E44').Value = '  -3.03%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '62.38'
$ws.Range('E45').Value = '  -0.26%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.711.19'
$ws.Range('E46').Value = '  -1.66%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '88.71'
$ws.Range('E47').Value = '  -2.04%  '

$ws.Range('E48').Value = '  -2.55%  '

$ws.Range('E49').Value = '  -4.42%  '

$ws.Range('E50').Value = '  -1.73%  '

$ws.Range('E51').Value = '  -1.66%  '
